$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sub = $s.Shapes.Item(2)
$sub.TextFrame.TextRange.Text = "2`r3`r4"
